# Weekly data update: a new price-report entry (date 2022-01-13, serial 44574)
# for "Vega Monumental Concepción" / Zanahoria is inserted right after the
# existing entry in row 97, pushing all subsequent rows down by two rows.
# The oldest entry that falls off the bottom (previously rows 173-174) ends
# up as the new last rows (175-176), so the net effect is a 2-row insert at
# row 98 with the new data filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 98; everything from old row 98 onward
# (up to old row 174) shifts down to rows 100-176.
$ws.Rows.Item(98).Resize(2).Insert()

# New row 98: "Primera" quality entry
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = "Vega Monumental Concepción"
$ws.Range("C98").Value = "Bíobío"
$ws.Range("D98").Value = 44574
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 100114013
$ws.Range("G98").Value = "Zanahoria"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8500
$ws.Range("M98").Value = 8250
$ws.Range("N98").Value = "`$/saco 20 kilos"
$ws.Range("O98").Value = "Región de Ñuble"
$ws.Range("P98").Value = 412
$ws.Range("Q98").Value = 20
$ws.Range("R98").Value = "Hortaliza"

# New row 99: "Segunda" quality entry
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 44574
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 100114013
$ws.Range("G99").Value = "Zanahoria"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Segunda"
$ws.Range("J99").Value = 500
$ws.Range("K99").Value = 7000
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = 7000
$ws.Range("N99").Value = "`$/saco 20 kilos"
$ws.Range("O99").Value = "Región de Ñuble"
$ws.Range("P99").Value = 350
$ws.Range("Q99").Value = 20
$ws.Range("R99").Value = "Hortaliza"
